$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Type"
$ws.Range("E1").Select() | Out-Null
